{"js": "// The document contains a date heading followed by a 5-column table of\n// two-digit-by-two-digit multiplication problems (\"NN\u00d7NN=\"). Every\n// non-empty paragraph in the body (the heading plus each populated table\n// cell, in document order) gets its text replaced with a new value.\n// We use an ordered list (old -> new) rather than a global text search and\n// replace because some source strings (e.g. \"54\u00d787=\") occur more than once\n// in the document but map to different replacement values depending on\n// their position.\n\nconst replacements = [\n  \"2024-05-14 Tuesday\",\n  \"23\u00d787=\",\n  \"57\u00d724=\",\n  \"36\u00d737=\",\n  \"94\u00d786=\",\n  \"48\u00d747=\",\n  \"22\u00d712=\",\n  \"97\u00d766=\",\n  \"16\u00d732=\",\n  \"83\u00d773=\",\n  \"62\u00d764=\",\n  \"18\u00d768=\",\n  \"69\u00d730=\",\n  \"98\u00d711=\",\n  \"22\u00d780=\",\n  \"79\u00d719=\",\n  \"29\u00d748=\",\n  \"23\u00d733=\",\n  \"33\u00d733=\",\n  \"28\u00d757=\",\n  \"63\u00d774=\",\n  \"74\u00d737=\",\n  \"14\u00d747=\",\n  \"41\u00d760=\",\n  \"86\u00d750=\",\n  \"93\u00d746=\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet replIndex = 0;\nfor (let i = 0; i < paragraphs.items.length && replIndex < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  // Skip the empty paragraphs that live in the blank table cells; only the\n  // paragraphs that actually hold one of the original \"NN\u00d7NN=\" / date\n  // strings should be touched.\n  if (para.text === \"\") {\n    continue;\n  }\n  para.insertText(replacements[replIndex], Word.InsertLocation.replace);\n  replIndex++;\n}\n\nawait context.sync();\n", "ps1": "# The document contains a date heading followed by a 5-column table of\n# two-digit-by-two-digit multiplication problems (\"NN\u00d7NN=\"). Every\n# non-empty paragraph in the body (the heading plus each populated table\n# cell, in document order) gets its text replaced with a new value.\n# We use an ordered list (old -> new) rather than a global Find/Replace\n# because some source strings (e.g. \"54\u00d787=\") occur more than once in the\n# document but map to different replacement values depending on their\n# position.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  \"2024-05-14 Tuesday\",\n  \"23\u00d787=\",\n  \"57\u00d724=\",\n  \"36\u00d737=\",\n  \"94\u00d786=\",\n  \"48\u00d747=\",\n  \"22\u00d712=\",\n  \"97\u00d766=\",\n  \"16\u00d732=\",\n  \"83\u00d773=\",\n  \"62\u00d764=\",\n  \"18\u00d768=\",\n  \"69\u00d730=\",\n  \"98\u00d711=\",\n  \"22\u00d780=\",\n  \"79\u00d719=\",\n  \"29\u00d748=\",\n  \"23\u00d733=\",\n  \"33\u00d733=\",\n  \"28\u00d757=\",\n  \"63\u00d774=\",\n  \"74\u00d737=\",\n  \"14\u00d747=\",\n  \"41\u00d760=\",\n  \"86\u00d750=\",\n  \"93\u00d746=\"\n)\n\n$idx = 0\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    if ($idx -ge $replacements.Count) { break }\n    $p = $d.Paragraphs($i)\n    $r = $p.Range\n    # Paragraphs that only contain the paragraph/cell end mark have\n    # End - Start == 1 (and blank table-cell paragraphs behave the same\n    # way); only touch paragraphs that actually hold text.\n    if (($r.End - $r.Start) -gt 1) {\n        $target = $d.Range($r.Start, $r.End - 1)\n        $target.Text = $replacements[$idx]\n        $idx += 1\n    }\n}\n"}
